# Apply crypto price/volume/coin updates (GitHub Actions refresh snapshot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cellRef, $text)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "62.999.11"
Set-TextCell "E2" "  +2.88%  "

Set-TextCell "D3" "2.455.44"
Set-TextCell "E3" "  +2.09%  "

Set-TextCell "E4" "  -0.14%  "

Set-TextCell "D5" "576.62"
Set-TextCell "E5" "  +1.45%  "

Set-TextCell "E6" "  +3.06%  "

Set-TextCell "E7" "  +0.12%  "

Set-TextCell "E8" "  +0.67%  "

Set-TextCell "D9" "2.455.29"
Set-TextCell "E9" "  +1.64%  "

Set-TextCell "D10" "0.111"
Set-TextCell "E10" "  +2.65%  "

Set-TextCell "E11" "  +2.41%  "

Set-TextCell "D12" "5.29"
Set-TextCell "E12" "  +1.13%  "

Set-TextCell "E13" "  +2.37%  "

Set-TextCell "D14" "28.55"
Set-TextCell "E14" "  +7.94%  "

Set-TextCell "D15" "0.0000179"
Set-TextCell "E15" "  +5.07%  "

Set-TextCell "D16" "2.898.68"
Set-TextCell "E16" "  +3.56%  "

Set-TextCell "D17" "62.945.08"
Set-TextCell "E17" "  +3.43%  "

Set-TextCell "D18" "2.470.55"
Set-TextCell "E18" "  +2.61%  "

Set-TextCell "D19" "7.94"
Set-TextCell "E19" "  -1.63%  "

Set-TextCell "D20" "11.08"
Set-TextCell "E20" "  +3.60%  "

Set-TextCell "D21" "330.09"
Set-TextCell "E21" "  +1.90%  "

Set-TextCell "D22" "4.13"
Set-TextCell "E22" "  +1.04%  "

Set-TextCell "E23" "  +10.09%  "

Set-TextCell "B25" "Litecoin"
Set-TextCell "C25" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D25" "66.42"
Set-TextCell "E25" "  +1.86%  "

Set-TextCell "B26" "Binance-PegBSC-USD"
Set-TextCell "C26" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextCell "D26" "1.25"
Set-TextCell "E26" "  +24.54%  "

Set-TextCell "D27" "658.94"
Set-TextCell "E27" "  +10.32%  "

Set-TextCell "D28" "8.65"
Set-TextCell "E28" "  +4.67%  "

Set-TextCell "E29" "  +6.37%  "

Set-TextCell "E30" "  +2.41%  "

Set-TextCell "D31" "8.23"
Set-TextCell "E31" "  +2.89%  "

Set-TextCell "E32" "  +4.66%  "

Set-TextCell "E33" "  +3.90%  "

Set-TextCell "E34" "  +4.85%  "

Set-TextCell "B35" "ImmutableX"
Set-TextCell "C35" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D35" "1.49"
Set-TextCell "E35" "  +1.61%  "

Set-TextCell "B36" "FirstDigitalUSD"
Set-TextCell "C36" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D36" "0.999"
Set-TextCell "E36" "  +0.13%  "

Set-TextCell "B37" "NEARProtocol"
Set-TextCell "C37" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D37" "4.80"
Set-TextCell "E37" "  +3.81%  "

Set-TextCell "B38" "RenderToken"
Set-TextCell "C38" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextCell "D38" "5.55"
Set-TextCell "E38" "  +4.93%  "

Set-TextCell "B39" "PolygonEcosystemToken"
Set-TextCell "C39" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextCell "D39" "0.374"
Set-TextCell "E39" "  +0.25%  "

Set-TextCell "B40" "EthereumClassic"
Set-TextCell "C40" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D40" "18.82"
Set-TextCell "E40" "  +2.67%  "

Set-TextCell "B41" "Monero"
Set-TextCell "C41" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D41" "151.70"
Set-TextCell "E41" "  +0.47%  "

Set-TextCell "B42" "dogwifhat"
Set-TextCell "C42" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell "D42" "2.74"
Set-TextCell "E42" "  +8.27%  "

Set-TextCell "B43" "Stacks"
Set-TextCell "C43" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D43" "1.77"
Set-TextCell "E43" "  +4.35%  "

Set-TextCell "B44" "OKB"
Set-TextCell "C44" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D44" "42.57"
Set-TextCell "E44" "  +1.34%  "

Set-TextCell "B45" "USDe"
Set-TextCell "C45" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell "D45" "1.00"
Set-TextCell "E45" "  +0.07%  "

Set-TextCell "B46" "WhiteBITCoin"
Set-TextCell "C46" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextCell "D46" "14.99"
Set-TextCell "E46" "  +27.16%  "

Set-TextCell "B47" "Aave"
Set-TextCell "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D47" "146.85"
Set-TextCell "E47" "  +3.75%  "

Set-TextCell "B48" "Filecoin"
Set-TextCell "C48" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D48" "3.63"
Set-TextCell "E48" "  +2.80%  "

Set-TextCell "B49" "InjectiveProtocol"
Set-TextCell "C49" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D49" "20.69"
Set-TextCell "E49" "  +3.98%  "

Set-TextCell "B50" "Mantle"
Set-TextCell "C50" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D50" "0.607"
Set-TextCell "E50" "  +2.36%  "

Set-TextCell "B51" "Hedera"
Set-TextCell "C51" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D51" "0.0516"
Set-TextCell "E51" "  +1.30%  "
